$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.380.21'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '3.509.75'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'591.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.69%  '
$ws.Range("D6").Value = "'134.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.32%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = "'7.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.94%  '
$ws.Range("E10").Value = '  +0.13%  '
$ws.Range("E11").Value = '  +3.20%  '
$ws.Range("D12").Value = '4.108.36'
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("D15").Value = '3.508.88'
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '64.362.70'
$ws.Range("E16").Value = '  -0.03%  '
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = "'25.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.60%  '
$ws.Range("D18").Value = "'10.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("E19").Value = '  +2.49%  '
$ws.Range("D20").Value = "'13.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.78%  '
$ws.Range("D21").Value = "'394.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.59%  '
$ws.Range("E22").Value = '  +1.21%  '
$ws.Range("D23").Value = '3.650.58'
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").Value = "'74.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.92%  '
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("D26").Value = "'5.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").Value = "'0.0000117"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.15%  '
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("D29").Value = "'7.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.86%  '
$ws.Range("D30").Value = "'2.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.03%  '
$ws.Range("D31").Value = "'8.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("E32").Value = '  -6.63%  '
$ws.Range("E33").Value = '  +5.71%  '
$ws.Range("D34").Value = '3.540.07'
$ws.Range("E34").Value = '  +0.42%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").Value = "'23.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.99%  '
$ws.Range("D37").Value = "'5.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.68%  '
$ws.Range("E38").Value = '  +1.20%  '
$ws.Range("E39").Value = '  +0.39%  '
$ws.Range("D40").Value = "'167.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("D41").Value = "'0.0786"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  +0.38%  '
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("D44").Value = "'25.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.91%  '
$ws.Range("E45").Value = '  +0.37%  '
$ws.Range("E46").Value = '  +0.62%  '
$ws.Range("E47").Value = '  -3.47%  '
$ws.Range("D48").Value = "'6.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.28%  '
$ws.Range("D49").Value = '2.381.01'
$ws.Range("E49").Value = '  -3.95%  '
$ws.Range("D50").Value = "'0.894"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.56%  '
$ws.Range("E51").Value = '  -0.12%  '
